# Update NI_vannf with the changes implemented in WFD2ECA
# Fill in min/max values for several water-quality parameters that
# previously had no numeric range (stored as "NA" text).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ANC - Syrenoytraliserende kapasitet (ANC)
$ws.Range("C16").Value = -200
$ws.Range("D16").Value = 250

# KLFA - Klorofyll a
$ws.Range("C157").Value = 0
$ws.Range("D157").Value = "Inf"

# LAL - Labilt aluminium
$ws.Range("C166").Value = 0
$ws.Range("D166").Value = "Inf"

# N-NH4 - Ammonium
$ws.Range("C258").Value = 0
$ws.Range("D258").Value = 1000000

# N-TOT - Totalnitrogen
$ws.Range("C270").Value = 0
$ws.Range("D270").Value = "Inf"

# PH - pH
$ws.Range("C278").Value = 3.2
$ws.Range("D278").Value = 10.8

# P-TOT - Totalfosfor
$ws.Range("C301").Value = 0
$ws.Range("D301").Value = "Inf"

# RAMI - River Acidification Macroinvertebrate Index (RAMI)
$ws.Range("C307").Value = 2
$ws.Range("D307").Value = 8

# SECCI - Siktedyp
$ws.Range("C315").Value = 0
$ws.Range("D315").Value = 60

# VANNSTVAR - Vannstandsvariasjoner (reguleringshoyde)
$ws.Range("C331").Value = 0
$ws.Range("D331").Value = 150
